# employee_import_template.xlsx — record a JSON-encoded allowances
# structure in the "allowances_json" column (AK) for both sample
# employee rows, instead of the old boolean flag values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AK column header ("allowances_json") now holds a JSON array describing
# a single "Housing" allowance of 15% instead of a plain TRUE/FALSE flag.
$allowancesJson = '[{"type": "Housing ", "value": 15, "calculation_type": "Percentage"}]'
$ws.Range("AK2").Value = $allowancesJson
$ws.Range("AK3").Value = $allowancesJson

# Widen column AK so the longer JSON sample value is readable.
$ws.Columns.Item(37).ColumnWidth = 49.45

# Selection / scroll position moved while reviewing the new column.
$ws.Range("AV2").Select()
$excel.ActiveWindow.ScrollColumn = 39
